# Workbook edit: update currency column (G) from MKD to EUR for loan rows 2-12,
# and move the active selection to G1 (also drops the stale topLeftCell scroll
# position from the previous session).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newCurrency = "ЕУР"

for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 7).Value = $newCurrency
}

# Update the active selection / view to match the saved workbook state.
$null = $ws.Range("G1").Select()
